$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text so that numeric-looking
# values (e.g. "0.1000", "0.830") keep their exact original formatting
# instead of being auto-converted into numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '57.204.28'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '2.422.28'
$ws.Range("E3").Value = '  -3.50%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.36%  '
$ws.Range("D5").Value = '489.33'
$ws.Range("E5").Value = '  -1.22%  '
$ws.Range("D6").Value = '155.12'
$ws.Range("E6").Value = '  +0.83%  '
$ws.Range("D7").Value = '0.617'
$ws.Range("E7").Value = '  +19.34%  '
$ws.Range("D8").Value = '0.997'
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '2.426.94'
$ws.Range("E9").Value = '  -4.00%  '
$ws.Range("D10").Value = '6.33'
$ws.Range("E10").Value = '  +9.93%  '
$ws.Range("D11").Value = '0.1000'
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("E12").Value = '  -1.56%  '
$ws.Range("E13").Value = '  +1.33%  '
$ws.Range("D14").Value = '2.839.90'
$ws.Range("E14").Value = '  -3.78%  '
$ws.Range("D15").Value = '57.172.01'
$ws.Range("E15").Value = '  -0.31%  '
$ws.Range("D16").Value = '20.62'
$ws.Range("E16").Value = '  -3.71%  '
$ws.Range("E17").Value = '  -3.72%  '
$ws.Range("D18").Value = '2.425.02'
$ws.Range("E18").Value = '  -3.82%  '
$ws.Range("E19").Value = '  +2.63%  '
$ws.Range("D20").Value = '324.87'
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("D21").Value = '10.02'
$ws.Range("E21").Value = '  -3.09%  '
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").Value = '57.88'
$ws.Range("E24").Value = '  -1.18%  '
$ws.Range("D25").Value = '0.405'
$ws.Range("E25").Value = '  -1.46%  '
$ws.Range("D26").Value = '0.997'
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("E27").Value = '  -2.08%  '
$ws.Range("D28").Value = '2.514.12'
$ws.Range("E28").Value = '  -3.98%  '
$ws.Range("D29").Value = '7.26'
$ws.Range("E29").Value = '  -4.80%  '
$ws.Range("D30").Value = '0.0₃0782'
$ws.Range("E30").Value = '  -6.45%  '
$ws.Range("E31").Value = '  +0.11%  '
$ws.Range("D32").Value = '151.08'
$ws.Range("E32").Value = '  -0.45%  '
$ws.Range("D33").Value = '18.61'
$ws.Range("E33").Value = '  +1.41%  '
$ws.Range("D34").Value = '1.53'
$ws.Range("E34").Value = '  -0.34%  '
$ws.Range("D35").Value = '5.29'
$ws.Range("E35").Value = '  -0.47%  '
$ws.Range("E36").Value = '  -1.19%  '
$ws.Range("D37").Value = '3.78'
$ws.Range("E37").Value = '  -1.77%  '
$ws.Range("D38").Value = '0.830'
$ws.Range("E38").Value = '  -6.99%  '
$ws.Range("D39").Value = '0.103'
$ws.Range("E39").Value = '  +8.94%  '
$ws.Range("D40").Value = '34.02'
$ws.Range("E40").Value = '  -0.94%  '
$ws.Range("E41").Value = '  -0.57%  '
$ws.Range("E42").Value = '  -3.32%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value = '279.75'
$ws.Range("E43").Value = '  +3.85%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '0.998'
$ws.Range("E44").Value = '  +0.33%  '
$ws.Range("E45").Value = '  -4.30%  '
$ws.Range("D46").Value = '0.0531'
$ws.Range("E46").Value = '  -5.81%  '
$ws.Range("D47").Value = '10.20'
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("E48").Value = '  -1.41%  '
$ws.Range("D49").Value = '4.52'
$ws.Range("E49").Value = '  -8.64%  '
$ws.Range("D50").Value = '1.897.60'
$ws.Range("E50").Value = '  -0.28%  '
$ws.Range("D51").Value = '17.56'
$ws.Range("E51").Value = '  -2.99%  '

# Restore the default (Normal) style on the price column so no stray
# number-format style is left behind on the cells.
$priceRange.Style = "Normal"
